# Auto-generated Excel COM-interop edit script
# Applies updated market-price / profit figures to the Kujata_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 2598.7
$ws.Range("I86").Value = 3666.6
$ws.Range("J86").Value = 1530.8
$ws.Range("K86").Value = 3666.6
$ws.Range("L86").Value = 1530.8
$ws.Range("M86").Value = -2543.6
$ws.Range("N86").Value = -3776.8

# Row 89
$ws.Range("H89").Value = 2598.7
$ws.Range("I89").Value = 3666.6
$ws.Range("J89").Value = 1530.8
$ws.Range("K89").Value = 18333
$ws.Range("L89").Value = 7654
$ws.Range("M89").Value = -12717
$ws.Range("N89").Value = -18886

# Row 137
$ws.Range("H137").Value = 2217.1155
$ws.Range("I137").Value = 1583
$ws.Range("J137").Value = 2901.96
$ws.Range("K137").Value = 4749
$ws.Range("L137").Value = 8705.880000000001
$ws.Range("M137").Value = -2199
$ws.Range("N137").Value = -13805.88

# Row 138
$ws.Range("H138").Value = 2750.2239
$ws.Range("I138").Value = 4325
$ws.Range("J138").Value = 2650.238
$ws.Range("K138").Value = 12975
$ws.Range("L138").Value = 7950.714
$ws.Range("M138").Value = -7835
$ws.Range("N138").Value = -18230.714

# Row 141
$ws.Range("H141").Value = 1266.6666
$ws.Range("I141").Value = 1099
$ws.Range("J141").Value = 2105
$ws.Range("K141").Value = 3297
$ws.Range("L141").Value = 6315
$ws.Range("M141").Value = 1883
$ws.Range("N141").Value = -16675

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1029.9
$ws.Range("I45").Value = 797.7143
$ws.Range("J45").Value = 1571.6666
$ws.Range("K45").Value = 797.7143
$ws.Range("L45").Value = 1571.6666
$ws.Range("M45").Value = -420.7143
$ws.Range("N45").Value = -2325.6666

# Row 61
$ws.Range("H61").Value = 125001976
$ws.Range("I61").Value = 166668380
$ws.Range("J61").Value = 2756.5
$ws.Range("K61").Value = 166668380
$ws.Range("L61").Value = 2756.5
$ws.Range("M61").Value = -166668168
$ws.Range("N61").Value = -3180.5

# Row 63
$ws.Range("H63").Value = 33335582
$ws.Range("I63").Value = 2296.1538
$ws.Range("J63").Value = 250001950
$ws.Range("K63").Value = 2296.1538
$ws.Range("L63").Value = 250001950
$ws.Range("M63").Value = -1610.1538
$ws.Range("N63").Value = -250003322

# Row 66
$ws.Range("H66").Value = 33335582
$ws.Range("I66").Value = 2296.1538
$ws.Range("J66").Value = 250001950
$ws.Range("K66").Value = 11480.769
$ws.Range("L66").Value = 1250009750
$ws.Range("M66").Value = -8048.769
$ws.Range("N66").Value = -1250016614

# Row 74
$ws.Range("H74").Value = 1291.2812
$ws.Range("I74").Value = 913.75
$ws.Range("J74").Value = 2423.875
$ws.Range("K74").Value = 913.75
$ws.Range("L74").Value = 2423.875
$ws.Range("M74").Value = -39.75
$ws.Range("N74").Value = -4171.875

# Row 77
$ws.Range("H77").Value = 1291.2812
$ws.Range("I77").Value = 913.75
$ws.Range("J77").Value = 2423.875
$ws.Range("K77").Value = 4568.75
$ws.Range("L77").Value = 12119.375
$ws.Range("M77").Value = -200.75
$ws.Range("N77").Value = -20855.375

# Row 136
$ws.Range("H136").Value = 125001976
$ws.Range("I136").Value = 166668380
$ws.Range("J136").Value = 2756.5
$ws.Range("K136").Value = 500005140
$ws.Range("L136").Value = 8269.5
$ws.Range("M136").Value = -500002590
$ws.Range("N136").Value = -13369.5

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 15168.25
$ws.Range("I82").Value = 3936.4285
$ws.Range("J82").Value = 30892.8
$ws.Range("K82").Value = 3936.4285
$ws.Range("L82").Value = 30892.8
$ws.Range("M82").Value = -3553.4285
$ws.Range("N82").Value = -31658.8

# Row 85
$ws.Range("H85").Value = 15168.25
$ws.Range("I85").Value = 3936.4285
$ws.Range("J85").Value = 30892.8
$ws.Range("K85").Value = 3936.4285
$ws.Range("L85").Value = 30892.8
$ws.Range("M85").Value = -2610.4285
$ws.Range("N85").Value = -33544.8

# Row 99
$ws.Range("H99").Value = 100001230
$ws.Range("I99").Value = 142858190
$ws.Range("J99").Value = 1650
$ws.Range("K99").Value = 142858190
$ws.Range("L99").Value = 1650
$ws.Range("M99").Value = -142856692
$ws.Range("N99").Value = -4646

# Row 105
$ws.Range("H105").Value = 126239440
$ws.Range("I105").Value = 144273360
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 144273360
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -144271613
$ws.Range("N105").Value = -5494

# Row 107
$ws.Range("H107").Value = 932.5
$ws.Range("I107").Value = 815.7143
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 815.7143
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 1104.2857
$ws.Range("N107").Value = -5590

# Row 134
$ws.Range("H134").Value = 5896.1904
$ws.Range("I134").Value = 1186.9375
$ws.Range("J134").Value = 20965.8
$ws.Range("K134").Value = 3560.8125
$ws.Range("L134").Value = 62897.39999999999
$ws.Range("M134").Value = -1025.8125
$ws.Range("N134").Value = -67967.39999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 8004.5
$ws.Range("I23").Value = 1009
$ws.Range("J23").Value = 15000
$ws.Range("K23").Value = 1009
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = -769
$ws.Range("N23").Value = -15480

# Row 27
$ws.Range("H27").Value = 8004.5
$ws.Range("I27").Value = 1009
$ws.Range("J27").Value = 15000
$ws.Range("K27").Value = 1009
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = -817
$ws.Range("N27").Value = -15384

# Row 31
$ws.Range("H31").Value = 1555.0193
$ws.Range("I31").Value = 1377.7906
$ws.Range("J31").Value = 2401.7778
$ws.Range("K31").Value = 1377.7906
$ws.Range("L31").Value = 2401.7778
$ws.Range("M31").Value = -1082.7906
$ws.Range("N31").Value = -2991.7778

# Row 34
$ws.Range("H34").Value = 1555.0193
$ws.Range("I34").Value = 1377.7906
$ws.Range("J34").Value = 2401.7778
$ws.Range("K34").Value = 1377.7906
$ws.Range("L34").Value = 2401.7778
$ws.Range("M34").Value = -1175.7906
$ws.Range("N34").Value = -2805.7778

# Row 35
$ws.Range("H35").Value = 675
$ws.Range("I35").Value = 675
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 675
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -381

# Row 58
$ws.Range("H58").Value = 5606.36
$ws.Range("I58").Value = 1164.1333
$ws.Range("J58").Value = 12269.7
$ws.Range("K58").Value = 1164.1333
$ws.Range("L58").Value = 12269.7
$ws.Range("M58").Value = -961.1333
$ws.Range("N58").Value = -12675.7

# Row 136
$ws.Range("H136").Value = 5606.36
$ws.Range("I136").Value = 1164.1333
$ws.Range("J136").Value = 12269.7
$ws.Range("K136").Value = 3492.3999
$ws.Range("L136").Value = 36809.10000000001
$ws.Range("M136").Value = -942.3998999999999
$ws.Range("N136").Value = -41909.10000000001

# Row 141
$ws.Range("H141").Value = 272812.34
$ws.Range("I141").Value = 2000
$ws.Range("J141").Value = 284586.78
$ws.Range("K141").Value = 2000
$ws.Range("L141").Value = 284586.78
$ws.Range("M141").Value = 3180
$ws.Range("N141").Value = -294946.78

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 33.625
$ws.Range("I2").Value = 26.666666
$ws.Range("J2").Value = 54.5
$ws.Range("K2").Value = 159.999996
$ws.Range("L2").Value = 327
$ws.Range("M2").Value = -46.99999600000001
$ws.Range("N2").Value = -553

# Row 4
$ws.Range("H4").Value = 4051215.2
$ws.Range("I4").Value = 2713324.8
$ws.Range("J4").Value = 4645833
$ws.Range("K4").Value = 8139974.399999999
$ws.Range("L4").Value = 13937499
$ws.Range("M4").Value = -8139862.399999999
$ws.Range("N4").Value = -13937723

# Row 6
$ws.Range("H6").Value = 51.333332
$ws.Range("I6").Value = 51.333332
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 153.999996
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -40.99999600000001

# Row 38
$ws.Range("H38").Value = 40
$ws.Range("I38").Value = 40
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 120
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 227
$ws.Range("N38").ClearContents()

# Row 138
$ws.Range("H138").Value = 3596.0476
$ws.Range("I138").Value = 4731.2856
$ws.Range("J138").Value = 3028.4285
$ws.Range("K138").Value = 14193.8568
$ws.Range("L138").Value = 9085.2855
$ws.Range("M138").Value = -9053.856800000001
$ws.Range("N138").Value = -19365.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 6895053
$ws.Range("I11").Value = 8000000
$ws.Range("J11").Value = 5000857
$ws.Range("K11").Value = 8000000
$ws.Range("L11").Value = 5000857
$ws.Range("M11").Value = -7999861
$ws.Range("N11").Value = -5001135

# Row 43
$ws.Range("H43").Value = 7000
$ws.Range("I43").Value = 7000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 7000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -6849

# Row 102
$ws.Range("H102").Value = 3779.75
$ws.Range("I102").Value = 2842
$ws.Range("J102").Value = 5342.6665
$ws.Range("K102").Value = 2842
$ws.Range("L102").Value = 5342.6665
$ws.Range("M102").Value = -1220
$ws.Range("N102").Value = -8586.666499999999

# Row 126
$ws.Range("H126").Value = 2275.5557
$ws.Range("I126").Value = 1756
$ws.Range("J126").Value = 2925
$ws.Range("K126").Value = 5268
$ws.Range("L126").Value = 8775
$ws.Range("M126").Value = -2798
$ws.Range("N126").Value = -13715

# Row 132
$ws.Range("H132").Value = 9534.611000000001
$ws.Range("I132").Value = 13581.2
$ws.Range("J132").Value = 4476.375
$ws.Range("K132").Value = 40743.60000000001
$ws.Range("L132").Value = 13429.125
$ws.Range("M132").Value = -38213.60000000001
$ws.Range("N132").Value = -18489.125

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2414.5386
$ws.Range("I7").Value = 2179.8333
$ws.Range("J7").Value = 2615.7144
$ws.Range("K7").Value = 2179.8333
$ws.Range("L7").Value = 2615.7144
$ws.Range("M7").Value = -2067.8333
$ws.Range("N7").Value = -2839.7144

# Row 31
$ws.Range("H31").Value = 3591.3333
$ws.Range("I31").Value = 262.5
$ws.Range("J31").Value = 4542.4287
$ws.Range("K31").Value = 262.5
$ws.Range("L31").Value = 4542.4287
$ws.Range("M31").Value = -14.5
$ws.Range("N31").Value = -5038.4287

# Row 32
$ws.Range("H32").Value = 8000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 8000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 8000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -8634

# Row 126
$ws.Range("H126").Value = 2414.5386
$ws.Range("I126").Value = 2179.8333
$ws.Range("J126").Value = 2615.7144
$ws.Range("K126").Value = 6539.499899999999
$ws.Range("L126").Value = 7847.1432
$ws.Range("M126").Value = -4069.499899999999
$ws.Range("N126").Value = -12787.1432

# Row 136
$ws.Range("H136").Value = 1688.2354
$ws.Range("I136").Value = 1600.8334
$ws.Range("J136").Value = 1898
$ws.Range("K136").Value = 4802.5002
$ws.Range("L136").Value = 5694
$ws.Range("M136").Value = -2252.5002
$ws.Range("N136").Value = -10794

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 1072.1578
$ws.Range("I136").Value = 710.6875
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2132.0625
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = 417.9375
$ws.Range("N136").Value = -14100
